# 2des pwbe aula07 add instrucoes
# Adds attendance marks ("P" = Presente / "F" = Falta) for a new lecture
# date in column AL (rows 3-30, row 12 is hidden and has no mark), and
# moves the active selection to AL31, mirroring the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AL3").Value  = "P"
$ws.Range("AL4").Value  = "P"
$ws.Range("AL5").Value  = "P"
$ws.Range("AL6").Value  = "P"
$ws.Range("AL7").Value  = "P"
$ws.Range("AL8").Value  = "P"
$ws.Range("AL9").Value  = "P"
$ws.Range("AL10").Value = "P"
$ws.Range("AL11").Value = "P"
# row 12 is hidden and is intentionally left untouched
$ws.Range("AL13").Value = "P"
$ws.Range("AL14").Value = "F"
$ws.Range("AL15").Value = "P"
$ws.Range("AL16").Value = "P"
$ws.Range("AL17").Value = "P"
$ws.Range("AL18").Value = "F"
$ws.Range("AL19").Value = "P"
$ws.Range("AL20").Value = "P"
$ws.Range("AL21").Value = "P"
$ws.Range("AL22").Value = "P"
$ws.Range("AL23").Value = "P"
$ws.Range("AL24").Value = "P"
$ws.Range("AL25").Value = "P"
$ws.Range("AL26").Value = "P"
$ws.Range("AL27").Value = "P"
$ws.Range("AL28").Value = "P"
$ws.Range("AL29").Value = "P"
$ws.Range("AL30").Value = "P"

# Move the selection to match where the author left off editing.
$ws.Range("AL31").Select()
